$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1628603333333333
$ws.Range("H2").Value = 0.488581
$ws.Range("I2").Value = 0.06904471801498467
$ws.Range("J2").Value = 0.06904471801498467
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.132884
$ws.Range("N2").Value = 0.398652
$ws.Range("O2").Value = 0.01195569974366677
$ws.Range("P2").Value = 0.01195569974366677
$ws.Range("Q2").Value = 0.02164153253466667
$ws.Range("R2").Value = 0.194773792812
$ws.Range("S2").Value = 0.0008254779174732964
$ws.Range("T2").Value = 0.0008254779174732964

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1628603333333333
$ws.Range("H3").Value = 0.488581
$ws.Range("I3").Value = 0.06904471801498467
$ws.Range("J3").Value = 0.06904471801498467
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07352966666666667
$ws.Range("N3").Value = 0.220589
$ws.Range("O3").Value = 0.006615533976389704
$ws.Range("P3").Value = 0.006615533976389703
$ws.Range("Q3").Value = 0.01197506602322222
$ws.Range("R3").Value = 0.107775594209
$ws.Range("S3").Value = 0.0004567676779183774
$ws.Range("T3").Value = 0.0004567676779183773

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1628603333333333
$ws.Range("H4").Value = 0.488581
$ws.Range("I4").Value = 0.06904471801498467
$ws.Range("J4").Value = 0.06904471801498467
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.908285
$ws.Range("N4").Value = 32.724855
$ws.Range("O4").Value = 0.9814287662799436
$ws.Range("P4").Value = 0.9814287662799435
$ws.Range("Q4").Value = 1.776526931195
$ws.Range("R4").Value = 15.988742380755
$ws.Range("S4").Value = 0.067762472419593
$ws.Range("T4").Value = 0.06776247241959299

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.769244333333333
$ws.Range("H5").Value = 5.307733
$ws.Range("I5").Value = 0.7500720009247772
$ws.Range("J5").Value = 0.7500720009247773
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.132884
$ws.Range("N5").Value = 0.398652
$ws.Range("O5").Value = 0.01195569974366677
$ws.Range("P5").Value = 0.01195569974366677
$ws.Range("Q5").Value = 0.2351042639906666
$ws.Range("R5").Value = 2.115938375916
$ws.Range("S5").Value = 0.008967635629187978
$ws.Range("T5").Value = 0.008967635629187978

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.769244333333333
$ws.Range("H6").Value = 5.307733
$ws.Range("I6").Value = 0.7500720009247772
$ws.Range("J6").Value = 0.7500720009247773
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.07352966666666667
$ws.Range("N6").Value = 0.220589
$ws.Range("O6").Value = 0.006615533976389704
$ws.Range("P6").Value = 0.006615533976389703
$ws.Range("Q6").Value = 0.1300919460818889
$ws.Range("R6").Value = 1.170827514737
$ws.Range("S6").Value = 0.004962126806856473
$ws.Range("T6").Value = 0.004962126806856473

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.769244333333333
$ws.Range("H7").Value = 5.307733
$ws.Range("I7").Value = 0.7500720009247772
$ws.Range("J7").Value = 0.7500720009247773
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.908285
$ws.Range("N7").Value = 32.724855
$ws.Range("O7").Value = 0.9814287662799436
$ws.Range("P7").Value = 0.9814287662799435
$ws.Range("Q7").Value = 19.299421422635
$ws.Range("R7").Value = 173.694792803715
$ws.Range("S7").Value = 0.7361422384887328
$ws.Range("T7").Value = 0.7361422384887328

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf16"
$ws.Range("C8").Value = "Fgfr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4266613333333333
$ws.Range("H8").Value = 1.279984
$ws.Range("I8").Value = 0.180883281060238
$ws.Range("J8").Value = 0.180883281060238
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.132884
$ws.Range("N8").Value = 0.398652
$ws.Range("O8").Value = 0.01195569974366677
$ws.Range("P8").Value = 0.01195569974366677
$ws.Range("Q8").Value = 0.05669646461866667
$ws.Range("R8").Value = 0.510268181568
$ws.Range("S8").Value = 0.002162586197005491
$ws.Range("T8").Value = 0.002162586197005491

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf16"
$ws.Range("C9").Value = "Fgfr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4266613333333333
$ws.Range("H9").Value = 1.279984
$ws.Range("I9").Value = 0.180883281060238
$ws.Range("J9").Value = 0.180883281060238
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.07352966666666667
$ws.Range("N9").Value = 0.220589
$ws.Range("O9").Value = 0.006615533976389704
$ws.Range("P9").Value = 0.006615533976389703
$ws.Range("Q9").Value = 0.03137226561955556
$ws.Range("R9").Value = 0.282350390576
$ws.Range("S9").Value = 0.001196639491614853
$ws.Range("T9").Value = 0.001196639491614853

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf16"
$ws.Range("C10").Value = "Fgfr4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4266613333333333
$ws.Range("H10").Value = 1.279984
$ws.Range("I10").Value = 0.180883281060238
$ws.Range("J10").Value = 0.180883281060238
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.908285
$ws.Range("N10").Value = 32.724855
$ws.Range("O10").Value = 0.9814287662799436
$ws.Range("P10").Value = 0.9814287662799435
$ws.Range("Q10").Value = 4.65414342248
$ws.Range("R10").Value = 41.88729080232
$ws.Range("S10").Value = 0.1775240553716176
$ws.Range("T10").Value = 0.1775240553716177
